$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.374.39"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "'1.716.36"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'225.05"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").Value = "'0.5278"
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("D7").Value = "'1.007"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.06668"
$ws.Range("E8").Value = "  +1.67%  "

$ws.Range("D9").Value = "'0.2652"
$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").Value = "'20.83"
$ws.Range("E10").Value = "  -1.25%  "

$ws.Range("D11").Value = "'0.07742"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "'4.475"
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("D13").Value = "'1.951.53"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").Value = "'1.716.56"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").Value = "'0.5802"
$ws.Range("E15").Value = "  +0.83%  "

$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").Value = "'67.85"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18").Value = "'27.369.66"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "'219.82"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").Value = "'1.008"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").Value = "'4.658"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").Value = "'10.43"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("D23").Value = "'6.052"
$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D25").Value = "'144.97"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").Value = "'0.1208"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").Value = "'7.230"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").Value = "'16.20"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("D30").Value = "'0.05350"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "'3.483"
$ws.Range("E32").Value = "  -0.69%  "

$ws.Range("D33").Value = "'3.399"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("D34").Value = "'1.642"
$ws.Range("E34").Value = "  +0.27%  "

$ws.Range("D35").Value = "'2.841"
$ws.Range("E35").Value = "  -1.17%  "

$ws.Range("D36").Value = "'0.9542"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").Value = "'2.400"
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("D38").Value = "'0.5884"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("D39").Value = "'1.185.72"
$ws.Range("E39").Value = "  +13.97%  "

$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").Value = "'0.8411"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").Value = "'101.15"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "'1.858.20"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").Value = "'0.0₈117"
$ws.Range("E46").Value = "  +2.07%  "

$ws.Range("D47").Value = "'57.58"
$ws.Range("E47").Value = "  -1.01%  "

$ws.Range("D48").Value = "'0.4548"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").Value = "'1.011"
$ws.Range("E49").Value = "  +0.67%  "

$ws.Range("D50").Value = "'8.168"
$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("E51").Value = "  -0.11%  "
